# Update countries & provincias Spain
# Applies the data refresh described in the commit: a handful of per-country
# case metrics change, the "Santa Lucia" / "Timor Oriental" rows swap names
# (their ranking order changed), and the "last updated" timestamp moves on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Countries whose ranking swapped position (Santa Lucia now ranks above
# --- Timor Oriental) -------------------------------------------------------
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Timor Oriental"

# --- Refreshed case counters ------------------------------------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 7339782
$ws.Range("C4").Value = 18439
$ws.Range("D4").Value = 4591057
$ws.Range("E4").Value = 2539120
$ws.Range("G4").Value = 152
$ws.Range("H4").Value = 209605

# Row 14: Francia
$ws.Range("D14").Value = 95426
$ws.Range("E14").Value = 415405

# Row 25: Alemania
$ws.Range("B25").Value = 288583
$ws.Range("C25").Value = 2245
$ws.Range("E25").Value = 28238
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 9545

# Row 73: Libano
$ws.Range("B73").Value = 37272
$ws.Range("C73").Value = 1018
$ws.Range("D73").Value = 16676
$ws.Range("E73").Value = 20245
$ws.Range("G73").Value = 4
$ws.Range("H73").Value = 351

# Row 116: Cabo Verde
$ws.Range("D116").Value = 5134
$ws.Range("E116").Value = 578
$ws.Range("G116").Value = 2
$ws.Range("H116").Value = 59

# Row 117: Malaui
$ws.Range("B117").Value = 5770
$ws.Range("C117").Value = 2
$ws.Range("D117").Value = 4243
$ws.Range("E117").Value = 1348

# Row 120: Suazilandia
$ws.Range("B120").Value = 5452
$ws.Range("C120").Value = 21
$ws.Range("D120").Value = 4844
$ws.Range("E120").Value = 500

# Row 128: Ruanda
$ws.Range("B128").Value = 4832
$ws.Range("C128").Value = 12
$ws.Range("D128").Value = 3117
$ws.Range("E128").Value = 1686

# Row 167: Republica del Chad
$ws.Range("B167").Value = 1185
$ws.Range("C167").Value = 7
$ws.Range("D167").Value = 1006
$ws.Range("E167").Value = 94
$ws.Range("G167").Value = 1
$ws.Range("H167").Value = 85

# --- Timestamp banner --------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 21:45"
